# Update destinations data with revised data.
# The data format had changed a bit so built that into the import file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 - Key Stage 4 (KS4) destinations: source link + refreshed date
$ws.Range("B11").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/find-statistics/key-stage-4-destination-measures/2022-23'>Key stage 4 destination measures</a>"
$ws.Range("C11").Value = "Aug 2022 -  Jul 2023 (21/22 learners) (27/02/25)"

# Row 12 - Key Stage 5 (KS5) destinations: source link + refreshed date
$ws.Range("B12").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/find-statistics/16-18-destination-measures'>16-18 destination measures</a>"
$ws.Range("C12").Value = "Aug 2022 -  Jul 2023 (21/22 learners) (27/02/25)"

# Reposition view/selection the way the author left it after the edit
$ws.Activate()
$ws.Range("B12").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
